# edit.ps1 - applies the Git.docx changes described by the commit diff:
#   1) Rewrite the "Git add -a" explanatory sentence (drop stray proofErr
#      grammar markers, fix "dot(.)" -> "dot (.)", capitalise "Instead",
#      and shrink "--a" -> "-a").
#   2) Move the "_GoBack" bookmark to sit between "Git " and "Push" in the
#      "Git Push origin master" paragraph (Word re-seats _GoBack itself -
#      adding it here removes the stale one near the end of the document).
#   3) Bold "Git checkout" and "-f" (leaving the separating space regular)
#      in the "Git checkout -f" paragraph.
#   4/5) Merge the split "git remote set-url --add --push origin " runs
#      (and drop their now-unnecessary spell-check proofErr wrapping) in
#      both "Remote Push URLs" examples.

$d = $word.ActiveDocument

# --- 1) "Git add -a" sentence rewrite -------------------------------------
$rng = $d.Content
$old1 = "staging area ( --a indicates to all data we can also use dot(.)instead of  --a"
$new1 = "staging area (--a indicates to all data we can also use dot (.)Instead of -a"
$rng.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)

# --- 2) Re-seat the _GoBack bookmark ---------------------------------------
$rng = $d.Content
$rng.Find.Execute("Git Push", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$gobackPos = $rng.Start + 4   # right after "Git " and before "Push"
$gobackRange = $d.Range($gobackPos, $gobackPos)
$d.Bookmarks.Add("_GoBack", $gobackRange)

# --- 3) Bold "Git checkout" and "-f" ---------------------------------------
$rng = $d.Content
$rng.Find.Execute("Git checkout", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Bold = 1

$rng = $d.Content
$rng.Find.Execute("–f", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Bold = 1

# --- 4/5) Merge "git remote set-url --add --push origin " runs -------------
$mergedUrlCmd = "git remote set-url --add --push origin "
$rng = $d.Content
$rng.Find.Execute($mergedUrlCmd, $true, $false, $false, $false, $false, $true, 1, $false, $mergedUrlCmd, 2)
# Continue the search on the same Range so it resumes after the first hit
# instead of re-matching the text we just rewrote.
$rng.Find.Execute($mergedUrlCmd, $true, $false, $false, $false, $false, $true, 1, $false, $mergedUrlCmd, 2)
